$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column G: "Stripe" credentials
$ws.Range("G1").Value = "Stripe"
$ws.Range("G2").Value = "womencoders@gmail.com"
$ws.Range("G3").Value = "Peoplespaceoc"

# Copy header style (bold) from F1 to G1
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

# Add hyperlink on G2 like E2, then reapply the built-in Hyperlink cell
# style by name so it reuses the existing style slot (matches E2's s="1")
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:womencoders@gmail.com")
$ws.Range("G2").Style = "Hyperlink"

# Set column width for G (matches F/E/D column sizing convention)
$ws.Range("G1").ColumnWidth = 23.736979166666668

# Update selection to match diff (G1 column selected)
$ws.Range("G1:G1048576").Select()
